# "update w/ fast data" - append newly logged weight readings to the
# raw_data log (Sheet2.xml / tab "raw_data"), rows 214:241, extending the
# existing date/time/weight/TOD table that currently ends at row 213.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# New readings (date serial, time-of-day fraction, weight) straight off the
# scale log, newest-first like the rest of the sheet.
$rows = @(
    @{Row=214; A=44108.425694444442; B=0.42569444444444443; C=71.3},
    @{Row=215; A=44108.425694444442; B=0.42569444444444443; C=71.3},
    @{Row=216; A=44108.425694444442; B=0.42569444444444443; C=71.3},
    @{Row=217; A=44108.367361111108; B=0.36736111111111108; C=71.3},
    @{Row=218; A=44108.366666666669; B=0.3666666666666667;  C=71.3},
    @{Row=219; A=44108.336805555555; B=0.33680555555555558; C=71.3},
    @{Row=220; A=44107.932638888888; B=0.93263888888888891; C=72.5},
    @{Row=221; A=44107.932638888888; B=0.93263888888888891; C=72.5},
    @{Row=222; A=44107.385416666664; B=0.38541666666666669; C=71.099999999999994},
    @{Row=223; A=44107.384722222225; B=0.38472222222222219; C=70.7},
    @{Row=224; A=44107.34097222222;  B=0.34097222222222223; C=70.7},
    @{Row=225; A=44107.340277777781; B=0.34027777777777773; C=70.7},
    @{Row=226; A=44106.881249999999; B=0.88124999999999998; C=70.7},
    @{Row=227; A=44110.317361111112; B=0.31736111111111115; C=70},
    @{Row=228; A=44110.316666666666; B=0.31666666666666665; C=70},
    @{Row=229; A=44110.272916666669; B=0.27291666666666664; C=69.400000000000006},
    @{Row=230; A=44110.148611111108; B=0.14861111111111111; C=69.900000000000006},
    @{Row=231; A=44109.936111111114; B=0.93611111111111101; C=69.900000000000006},
    @{Row=232; A=44109.88958333333;  B=0.88958333333333339; C=69.900000000000006},
    @{Row=233; A=44109.729166666664; B=0.72916666666666663; C=70.400000000000006},
    @{Row=234; A=44109.34097222222;  B=0.34097222222222223; C=72.3},
    @{Row=235; A=44109.34097222222;  B=0.34097222222222223; C=72},
    @{Row=236; A=44109.306250000001; B=0.30624999999999997; C=72},
    @{Row=237; A=44109.306250000001; B=0.30624999999999997; C=72},
    @{Row=238; A=44109.305555555555; B=0.30555555555555552; C=72.599999999999994},
    @{Row=239; A=44109.292361111111; B=0.29236111111111113; C=71.7},
    @{Row=240; A=44108.915972222225; B=0.9159722222222223;  C=72.5},
    @{Row=241; A=44108.915277777778; B=0.91527777777777775; C=72.900000000000006}
)

$lastRow = 213

foreach ($r in $rows) {
    $row = $r.Row

    # Pull the date/time number formats (and row styling) down from the
    # last existing data row so the new rows look identical to the old ones.
    $ws.Range("A$lastRow`:D$lastRow").Copy() | Out-Null
    $ws.Range("A$row`:D$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Formula = "=IF(B$row<TIME(12,0,0), ""AM"", ""PM"")"

    $lastRow = $row
}

# Match the post-edit selection/scroll state recorded in the workbook.
$ws.Range("A242").Select()
